# Commit: "case with 380 kV done"
# Recomputed res_line/loading_percent.xlsx results for Case_2_232 after the
# 380 kV network change: overwrite the per-line loading-percentage values
# (columns B,C,E,F,G,H,J,K,M,N,O; rows 2-25) with their newly simulated
# values, leaving the hour index (col A) and the always-zero columns
# (D,I,L) untouched.
# Data-driven: each line below is "Row,Col,Value" describing one new cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function ColToNum($col) {
    $n = 0
    foreach ($c in $col.ToCharArray()) {
        $n = $n * 26 + ([int][char]$c - [int][char]'A' + 1)
    }
    return $n
}

$cellData = @"
2,B,8.017391377862316
2,C,6.018113802680668
2,E,22.12849206409117
2,F,38.85650716094047
2,G,25.8361480638682
2,H,13.67937791916493
2,J,7.89530771654847
2,K,7.878620275915256
2,M,17.66896689076781
2,N,18.40222474829963
2,O,20.40989766955055
3,B,7.738672606967256
3,C,5.958268909931895
3,E,22.02530026217837
3,F,38.79845731967963
3,G,25.93025796691386
3,H,13.72399933616677
3,J,7.905200320535685
3,K,7.642710755869878
3,M,17.54768505521361
3,N,18.45518949955821
3,O,20.48662678492748
4,B,7.563355172737432
4,C,5.921061016371024
4,E,21.96641309066864
4,F,38.77233371705716
4,G,25.99568874018788
4,H,13.75326303998146
4,J,7.911602856321928
4,K,7.492636344383291
4,M,17.47585520772737
4,N,18.48930729277505
4,O,20.53756806936188
5,B,7.49097121952632
5,C,5.905791312424662
5,E,21.94355919066978
5,F,38.76408988887717
5,G,26.02426733079965
5,H,13.76565784100288
5,J,7.914294748497524
5,K,7.430221016511041
5,M,17.44727204810579
5,N,18.50361329371662
5,O,20.55928911354027
6,B,7.478898553996609
6,C,5.903249586429136
6,E,21.93983386333273
6,F,38.76286626185731
6,G,26.029128209857
6,H,13.76774436284918
6,J,7.914746742910086
6,K,7.419782680421472
6,M,17.44256809843448
6,N,18.50601314986881
6,O,20.56295395745728
7,B,7.562382631784273
7,C,5.92085550557362
7,E,21.96610022357565
7,F,38.77221280423993
7,G,25.99606641787197
7,H,13.7534282986717
7,J,7.911638824548007
7,K,7.491799609733355
7,M,17.47546690725156
7,N,18.48949859644975
7,O,20.53785711275462
8,B,7.922220862343685
8,C,5.997580856260671
8,E,22.0919946148232
8,F,38.83452131629334
8,G,25.86700542257808
8,H,13.69437633508443
8,J,7.898650651331374
8,K,7.798394611524559
8,M,17.62661667632374
8,N,18.42015614844649
8,O,20.4355585889906
9,B,8.590304149427935
9,C,6.143967430550576
9,E,22.37348986265347
9,F,39.03179384262916
9,G,25.67492785716806
9,H,13.59336365980945
9,J,7.875776852388992
9,K,8.356002975074814
9,M,17.94281137721896
9,N,18.29679877473299
9,O,20.26537306064844
10,B,9.053098292261037
10,C,6.248511427078
10,E,22.60014149472659
10,F,39.22177606864712
10,G,25.57143735408885
10,H,13.52814044602749
10,J,7.860539934060638
10,K,8.736524303800126
10,M,18.18559917245002
10,N,18.21379217195981
10,O,20.15893492894987
11,B,9.256645670424103
11,C,6.295310217771699
11,E,22.70726179299561
11,F,39.31778645921385
11,G,25.53261181363665
11,H,13.50041602150896
11,J,7.853945840204402
11,K,8.902851886639748
11,M,18.29797286081721
11,N,18.17767069429978
11,O,20.11456277381835
12,B,9.332656691550858
12,C,6.312913710258193
12,E,22.74837607882073
12,F,39.35550148523335
12,G,25.5191024827933
12,H,13.49019700131206
12,J,7.851497100835616
12,K,8.964829537150282
12,M,18.34077272504161
12,N,18.16422695668088
12,O,20.09834314268864
13,B,9.316334795012974
13,C,6.309127882033786
13,E,22.73949735766775
13,F,39.34731883175726
13,G,25.52195880243613
13,H,13.49238541724304
13,J,7.852022335585357
13,K,8.951526812623458
13,M,18.33154454550885
13,N,18.16711188360416
13,O,20.10181037645714
14,B,9.262920886219053
14,C,6.296760887672502
14,E,22.71063342366799
14,F,39.32086220609269
14,G,25.53147644791099
14,H,13.4995696937526
14,J,7.85374341403945
14,K,8.907971182408755
14,M,18.30148927150016
14,N,18.176559971949
14,O,20.11321668127896
15,B,9.230062474060398
15,C,6.289170088599332
15,E,22.69302425152082
15,F,39.30483296320515
15,G,25.53746183036872
15,H,13.50400667910741
15,J,7.854803908691334
15,K,8.881160046637433
15,M,18.2831107129146
15,N,18.18237773333923
15,O,20.12027934928463
16,B,9.039649863498601
16,C,6.245436965838983
16,E,22.59321942385352
16,F,39.21569275790959
16,G,25.57414131395871
16,H,13.52999144178389
16,J,7.860977643254855
16,K,8.725515199141491
16,M,18.17829145607274
16,N,18.21618568741044
16,O,20.16191629191468
17,B,8.921001887970629
17,C,6.218407234967022
17,E,22.5330017551254
17,F,39.16344948315381
17,G,25.59876139125848
17,H,13.54643050387247
17,J,7.864851268800339
17,K,8.62827259003598
17,M,18.11446024593796
17,N,18.2373447858331
17,O,20.188496600925
18,B,8.852105029324342
18,C,6.202789556730419
18,E,22.49874578546133
18,F,39.13430422223232
18,G,25.61369865683387
18,H,13.55606899047998
18,J,7.867111031582062
18,K,8.571705379638113
18,M,18.07792959787852
18,N,18.24966923250164
18,O,20.20416575241245
19,B,8.828667555094921
19,C,6.197489783757612
19,E,22.48721332005246
19,F,39.12459194047659
19,G,25.61888931302155
19,H,13.55936388305082
19,J,7.86788160853282
19,K,8.552444556804293
19,M,18.06559337908753
19,N,18.25386860758174
19,O,20.20953643519668
20,B,8.933700311601482
20,C,6.221291999496697
20,E,22.53937293156554
20,F,39.16891748303053
20,G,25.59606013508018
20,H,13.54466158147202
20,J,7.864435629566849
20,K,8.63869028911734
20,M,18.12123643132334
20,N,18.23507640188576
20,O,20.18562765545663
21,B,9.278639302862111
21,C,6.300396652677328
21,E,22.71909675171928
21,F,39.32859647702384
21,G,25.52864845777346
21,H,13.4974519104353
21,J,7.853236582261816
21,K,8.920792100467406
21,M,18.31031079489712
21,N,18.1737784777032
21,O,20.10985053744402
22,B,9.497827646081113
22,C,6.3514031450448
22,E,22.83975007562245
22,F,39.44086032107748
22,G,25.49154763596107
22,H,13.46822739974464
22,J,7.846198802010735
22,K,9.099279709904991
22,M,18.4353034400673
22,N,18.135084214216
22,O,20.06372534313076
23,B,9.381433724410375
23,C,6.324246366934956
23,E,22.77507217501672
23,F,39.38022704469946
23,G,25.51071056578107
23,H,13.48367601505828
23,J,7.849929309976952
23,K,9.004565670439726
23,M,18.36847271403654
23,N,18.1556112584147
23,O,20.08803178328552
24,B,8.927961481556041
24,C,6.219988039495855
24,E,22.53649138688628
24,F,39.16644262615565
24,G,25.59727893494052
24,H,13.54546072737108
24,J,7.864623437992887
24,K,8.633982504145816
24,M,18.11817239609111
24,N,18.23610144042452
24,O,20.18692349788251
25,B,8.414160016510756
25,C,6.104856518901243
25,E,22.29375730743153
25,F,38.97045470074751
25,G,25.72031168885837
25,H,13.61910967085952
25,J,7.875776852388992
25,K,8.21010745760629
25,M,17.85531416824893
25,N,20.26537306064844
25,O,20.259901154675234
"@

$lines = $cellData -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split ","
    $row = [int]$parts[0]
    $col = $parts[1]
    $value = [double]$parts[2]
    $colNum = ColToNum $col
    $ws.Cells.Item($row, $colNum).Value = $value
}
